# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across several
# Leve-profit worksheets, reflecting refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row64: "Forged from the Void"
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2800
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 2800
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -3296

# ALC!row67: "Dodging the Draft (L)"
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 2800
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -4516

# ALC!row76: "Warding Off Temptation"
$ws.Range("H76").Value = 4200
$ws.Range("I76").Value = 3400
$ws.Range("K76").Value = 3400
$ws.Range("M76").Value = -3085

# ALC!row79: "The Garden of Arcane Delights (L)"
$ws.Range("H79").Value = 4200
$ws.Range("I79").Value = 3400
$ws.Range("K79").Value = 3400
$ws.Range("M79").Value = -2308

# ALC!row132: "Fast-forwarding Flora"
$ws.Range("H132").Value = 16735076
$ws.Range("I132").Value = 18944792
$ws.Range("J132").Value = 4363.7144
$ws.Range("K132").Value = 56834376
$ws.Range("L132").Value = 13091.1432
$ws.Range("M132").Value = -56831846
$ws.Range("N132").Value = -18151.1432

$ws = $wb.Worksheets.Item("ARM")
# ARM!row2: "Ain't Got No Ingots"
$ws.Range("H2").Value = 786.5
$ws.Range("I2").Value = 747.2273
$ws.Range("K2").Value = 747.2273
$ws.Range("M2").Value = -634.2273

# ARM!row64: "Don't Scuttle with Scuta"
$ws.Range("H64").Value = 48876.43
$ws.Range("J64").Value = 48876.43
$ws.Range("L64").Value = 48876.43
$ws.Range("N64").Value = -49372.43

# ARM!row67: "Shielded by Bureaucracy (L)"
$ws.Range("H67").Value = 48876.43
$ws.Range("J67").Value = 48876.43
$ws.Range("L67").Value = 48876.43
$ws.Range("N67").Value = -50592.43

# ARM!row68: "Let Faith Light the Way"
$ws.Range("H68").Value = 50099
$ws.Range("J68").Value = 50099
$ws.Range("L68").Value = 50099
$ws.Range("N68").Value = -51721

# ARM!row71: "Fifty Shields of Blades (L)"
$ws.Range("H71").Value = 50099
$ws.Range("J71").Value = 50099
$ws.Range("L71").Value = 150297
$ws.Range("N71").Value = -158409

# ARM!row80: "A Squire to Inspire"
$ws.Range("H80").Value = 36124.777
$ws.Range("J80").Value = 36124.777
$ws.Range("L80").Value = 36124.777
$ws.Range("N80").Value = -38120.777

# ARM!row83: "All's Fair in Highborn Assassination (L)"
$ws.Range("H83").Value = 36124.777
$ws.Range("J83").Value = 36124.777
$ws.Range("L83").Value = 108374.331
$ws.Range("N83").Value = -118358.331

# ARM!row88: "The Mast Chance"
$ws.Range("H88").Value = 8336633
$ws.Range("J88").Value = 3800
$ws.Range("L88").Value = 3800
$ws.Range("N88").Value = -4612

# ARM!row91: "The Rose and the Riveter (L)"
$ws.Range("H91").Value = 8336633
$ws.Range("J91").Value = 3800
$ws.Range("L91").Value = 3800
$ws.Range("N91").Value = -6608

# ARM!row116: "No Scope"
$ws.Range("H116").Value = 786.5
$ws.Range("I116").Value = 747.2273
$ws.Range("K116").Value = 747.2273
$ws.Range("M116").Value = 1546.7727

# ARM!row132: "Don't Bore Me, Ore Me"
$ws.Range("H132").Value = 1785.2452
$ws.Range("I132").Value = 841.0789
$ws.Range("K132").Value = 2523.2367
$ws.Range("M132").Value = 6.763300000000072

$ws = $wb.Worksheets.Item("BSM")
# BSM!row3: "Hells Bells"
$ws.Range("H3").Value = 786.5
$ws.Range("I3").Value = 747.2273
$ws.Range("K3").Value = 747.2273
$ws.Range("M3").Value = -633.2273

# BSM!row86: "Through Thick and Thin"
$ws.Range("H86").Value = 2272
$ws.Range("I86").Value = 2223.1428
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 2223.1428
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -1100.1428
$ws.Range("N86").Value = -4746

# BSM!row89: "Piercing Eyes Deserve Piercing Shafts (L)"
$ws.Range("H89").Value = 2272
$ws.Range("I89").Value = 2223.1428
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 11115.714
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -5499.714
$ws.Range("N89").Value = -23732

# BSM!row105: "Ingot to Wing It"
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = $null

$ws = $wb.Worksheets.Item("CRP")
# CRP!row31: "Wall Not Found"
$ws.Range("H31").Value = 2617.7727
$ws.Range("I31").Value = 1000.70966
$ws.Range("K31").Value = 1000.70966
$ws.Range("M31").Value = -705.70966

# CRP!row34: "Armoires of the Rich and Famous"
$ws.Range("H34").Value = 2617.7727
$ws.Range("I34").Value = 1000.70966
$ws.Range("K34").Value = 1000.70966
$ws.Range("M34").Value = -798.70966

# CRP!row62: "Splinter in the Sewers"
$ws.Range("H62").Value = 100004400
$ws.Range("I62").Value = 125004250
$ws.Range("J62").Value = 5006
$ws.Range("K62").Value = 125004250
$ws.Range("L62").Value = 5006
$ws.Range("M62").Value = -125003626
$ws.Range("N62").Value = -6254

# CRP!row65: "The Lumber of Their Discontent (L)"
$ws.Range("H65").Value = 100004400
$ws.Range("I65").Value = 125004250
$ws.Range("J65").Value = 5006
$ws.Range("K65").Value = 625021250
$ws.Range("L65").Value = 25030
$ws.Range("M65").Value = -625018130
$ws.Range("N65").Value = -31270

# CRP!row86: "Birch, Please"
$ws.Range("H86").Value = 3066.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3066.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3066.5
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -5312.5

# CRP!row89: "Built This City on Blocks and Soul (L)"
$ws.Range("H89").Value = 3066.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3066.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 15332.5
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -26564.5

# CRP!row97: "Wood That You Could"
$ws.Range("H97").Value = 34373.75
$ws.Range("J97").Value = 34373.75
$ws.Range("L97").Value = 34373.75
$ws.Range("N97").Value = -36355.75

# CRP!row134: "Wood You Be Quiet"
$ws.Range("H134").Value = 3583
$ws.Range("I134").Value = 3740.1177
$ws.Range("J134").Value = 3172.077
$ws.Range("K134").Value = 11220.3531
$ws.Range("L134").Value = 9516.231
$ws.Range("M134").Value = -8685.3531
$ws.Range("N134").Value = -14586.231

$ws = $wb.Worksheets.Item("GSM")
# GSM!row70: "Sky Is the Limit"
$ws.Range("H70").Value = 6223.2144
$ws.Range("I70").Value = 5896.1724
$ws.Range("K70").Value = 5896.1724
$ws.Range("M70").Value = -5626.1724

# GSM!row73: "Hulls of Broken Dreams (L)"
$ws.Range("H73").Value = 6223.2144
$ws.Range("I73").Value = 5896.1724
$ws.Range("K73").Value = 5896.1724
$ws.Range("M73").Value = -4960.1724

# GSM!row80: "Needs More Prayerbell"
$ws.Range("H80").Value = 16669460
$ws.Range("I80").Value = 41668884
$ws.Range("K80").Value = 41668884
$ws.Range("M80").Value = -41667886

# GSM!row83: "With a Noise That Reaches Heaven (L)"
$ws.Range("H83").Value = 16669460
$ws.Range("I83").Value = 41668884
$ws.Range("K83").Value = 208344420
$ws.Range("M83").Value = -208339428

# GSM!row132: "On Board for Lar"
$ws.Range("H132").Value = 2508.7256
$ws.Range("I132").Value = 1690.8125
$ws.Range("K132").Value = 5072.4375
$ws.Range("M132").Value = -2542.4375
